$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns H:L entirely (shrinking the used range from A1:L9 to A1:G9)
$ws.Range("H1:L9").Delete()

# Row 1: years
$ws.Range("B1").Value = 2015
$ws.Range("C1").Value = 2016
$ws.Range("D1").Value = 2017
$ws.Range("E1").Value = 2018
$ws.Range("F1").Value = 2019
$ws.Range("G1").Value = 2020

# Row 2: count
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 42
$ws.Range("D2").Value = 42
$ws.Range("E2").Value = 42
$ws.Range("F2").Value = 41
$ws.Range("G2").Value = 44

# Row 3: mean
$ws.Range("C3").Value = 5.723558153020192
$ws.Range("D3").Value = 5.798152265637031
$ws.Range("E3").Value = 6.117503522452244
$ws.Range("F3").Value = 5.827453470229446
$ws.Range("G3").Value = 5.561777815372469

# Row 4: std
$ws.Range("C4").Value = 1.253763836537182
$ws.Range("D4").Value = 1.183294736722819
$ws.Range("E4").Value = 1.167959067398114
$ws.Range("F4").Value = 0.9344294114655046
$ws.Range("G4").Value = 0.8808151154418762

# Row 5: min
$ws.Range("C5").Value = 3.152336191851997
$ws.Range("D5").Value = 3.628306484795113
$ws.Range("E5").Value = 3.322370620528491
$ws.Range("F5").Value = 4.220876693802421
$ws.Range("G5").Value = 4.36

# Row 6: 25%
$ws.Range("C6").Value = 4.864806790115457
$ws.Range("D6").Value = 4.958984573235083
$ws.Range("E6").Value = 5.607525680470014
$ws.Range("F6").Value = 5.155529335332706
$ws.Range("G6").Value = 4.798681519314282

# Row 7: 50%
$ws.Range("C7").Value = 5.674900271730575
$ws.Range("D7").Value = 5.948942443930747
$ws.Range("E7").Value = 6.277415476188198
$ws.Range("F7").Value = 5.838999475320748
$ws.Range("G7").Value = 5.4581474413627

# Row 8: 75%
$ws.Range("C8").Value = 6.859348464173619
$ws.Range("D8").Value = 6.900025957304678
$ws.Range("E8").Value = 6.827418096368399
$ws.Range("F8").Value = 6.219868229362246
$ws.Range("G8").Value = 6.269645442077534

# Row 9: max
$ws.Range("C9").Value = 8.102198711914406
$ws.Range("D9").Value = 7.563904173416154
$ws.Range("E9").Value = 7.895786943097561
$ws.Range("F9").Value = 7.754257368628767
$ws.Range("G9").Value = 7.232001944008197
